# adicionando tags IPI e iCMS (+ tabela de-para)
#
# Insere duas novas colunas na planilha:
#   - "DESC CST ICMS" (nova coluna T) : descricao do CST ICMS (coluna S),
#     obtida por uma tabela de-para codigo -> descricao.
#   - "TIPI" (nova coluna V)          : nova coluna de controle, preenchida
#     com "NAO" para as notas existentes.
#
# As colunas antigas T..Y (IPI_CST, CONFINS, Sujeito a ISS?, Outros Impostos,
# Infos Adicionais, DIFAL) sao deslocadas para a direita automaticamente
# pelo Insert das novas colunas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insere da direita para a esquerda para que os enderecos das colunas ainda
# nao processadas permanecam validos.

# 1) Nova coluna logo antes da antiga coluna U ("CONFINS") -> vira "TIPI"
$ws.Columns("U:U").Insert()

# 2) Nova coluna logo antes da antiga coluna T ("IPI_CST") -> vira "DESC CST ICMS"
$ws.Columns("T:T").Insert()

# Cabecalhos das novas colunas
$ws.Range("T1").Value = "DESC CST ICMS"
$ws.Range("V1").Value = "TIPI"

# Tabela de-para: codigo CST ICMS (coluna S) -> descricao (nova coluna "DESC CST ICMS")
$tabelaDeParaCstIcms = @{
    "00" = "Tributada integralmente"
    "15" = "CST 15"
    "20" = "Com redução da base de cálculo"
    "53" = "CST 53"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, "S").End(-4162).Row

for ($row = 2; $row -le $lastRow; $row++) {
    $cstIcms = $ws.Range("S$row").Value()
    if ($tabelaDeParaCstIcms.ContainsKey($cstIcms)) {
        $descricao = $tabelaDeParaCstIcms[$cstIcms]
    } else {
        $descricao = "CST $cstIcms"
    }
    $ws.Range("T$row").Value = $descricao
    $ws.Range("V$row").Value = "NAO"
}
